$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ramp start time (E1) and the standalone timepoint E4.
$ws.Range("E1").Value = 0.54583333333333328
$ws.Range("E4").Value = 0.58888888888888891

# Rows 16 and 19 no longer follow the incremental formula chain - they now
# hold directly observed/typed time values, breaking the shared formula.
$ws.Range("E16").Value = 0.76250000000000007
$ws.Range("E19").Value = 0.8041666666666667

# Recalculate so the remaining formula cells (E2,E3,E5:E15,E17,E18,E20:E24)
# pick up the new upstream values.
$excel.Calculate()

# Update the active selection to match the saved view.
$ws.Range("G19").Select()
